# Correction des tests suite a la revision de Boivin
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections de notes (colonne E "C4" et colonne H "Note sur 100") ---
# Plusieurs eleves avaient des notes erronees ramenees a 0 apres revision.
$ws.Range("E11").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("E37").Value = 0

# --- Nouvelle colonne "Somme" (J) = somme des criteres C1:C6 (B:G) ---
$ws.Range("J2").Value = "Somme"
$ws.Range("J3").Formula = "=SUM(B3:G3)"
$ws.Range("J4:J37").Formula = "=SUM(B4:G4)"

# --- Correlation entre la somme et la note sur 100 ---
$ws.Range("L4").Formula = "=CORREL(J3:J37,H3:H37)"

# --- Selection / vue active ---
[void]$ws.Range("E38").Select()
